# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (blank) column between the
# existing "In Advance" (M) and "Late" (N) columns, pushing "Late",
# "Outstanding" heading and the trailing "Outstanding" values each one
# column to the right (N->O, O->P, P->Q). The new column inherits its
# width from the column immediately to its left ("In Advance").
#
# The previously-active tab ("Transactions") is no longer the selected
# sheet; "Repayment schedule" becomes the active sheet/tab instead, with
# its selection parked at K19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting
# N/O/P -> O/P/Q and carrying over the formatting of column N.
$ws.Columns("N:N").Insert()

# Match the new column's width to the column to its left ("In Advance").
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and park the selection
# at K19 (was previously on the "Transactions" sheet at D3).
$ws.Activate()
$ws.Range("K19").Select()
